# Auto-generated script applying market-data refresh values to the "Leve Profit" sheets.
# Each entry updates one cell (H..N) on a given sheet/row to match the refreshed
# market snapshot. A $null Value means the cell is fully cleared (removed), matching
# rows where no market data is available for that column.

$wb = $excel.ActiveWorkbook

$changes = @(
    # --- ALC ---
    @{Sheet="ALC"; Cell="H19"; Value=1352.9333},
    @{Sheet="ALC"; Cell="I19"; Value=903.6667},
    @{Sheet="ALC"; Cell="K19"; Value=903.6667},
    @{Sheet="ALC"; Cell="M19"; Value=-728.6667},
    @{Sheet="ALC"; Cell="H21"; Value=8666.666999999999},
    @{Sheet="ALC"; Cell="I21"; Value=2000},
    @{Sheet="ALC"; Cell="J21"; Value=12000},
    @{Sheet="ALC"; Cell="K21"; Value=2000},
    @{Sheet="ALC"; Cell="L21"; Value=12000},
    @{Sheet="ALC"; Cell="M21"; Value=-1532},
    @{Sheet="ALC"; Cell="N21"; Value=-12936},
    @{Sheet="ALC"; Cell="H23"; Value=8666.666999999999},
    @{Sheet="ALC"; Cell="I23"; Value=2000},
    @{Sheet="ALC"; Cell="J23"; Value=12000},
    @{Sheet="ALC"; Cell="K23"; Value=2000},
    @{Sheet="ALC"; Cell="L23"; Value=12000},
    @{Sheet="ALC"; Cell="M23"; Value=-1766},
    @{Sheet="ALC"; Cell="N23"; Value=-12468},
    @{Sheet="ALC"; Cell="H33"; Value=879.2727},
    @{Sheet="ALC"; Cell="I33"; Value=185.33333},
    @{Sheet="ALC"; Cell="K33"; Value=185.33333},
    @{Sheet="ALC"; Cell="M33"; Value=43.66667000000001},
    @{Sheet="ALC"; Cell="H34"; Value=1500},
    @{Sheet="ALC"; Cell="I34"; Value=1500},
    @{Sheet="ALC"; Cell="K34"; Value=1500},
    @{Sheet="ALC"; Cell="M34"; Value=-1297},
    @{Sheet="ALC"; Cell="H36"; Value=1500},
    @{Sheet="ALC"; Cell="I36"; Value=1500},
    @{Sheet="ALC"; Cell="K36"; Value=1500},
    @{Sheet="ALC"; Cell="M36"; Value=-785},
    @{Sheet="ALC"; Cell="H40"; Value=3205.3157},
    @{Sheet="ALC"; Cell="J40"; Value=4375.1113},
    @{Sheet="ALC"; Cell="L40"; Value=4375.1113},
    @{Sheet="ALC"; Cell="N40"; Value=-4725.1113},
    @{Sheet="ALC"; Cell="H112"; Value=110984.9},
    @{Sheet="ALC"; Cell="I112"; Value=1231.25},
    @{Sheet="ALC"; Cell="J112"; Value=549999.5},
    @{Sheet="ALC"; Cell="K112"; Value=3693.75},
    @{Sheet="ALC"; Cell="L112"; Value=1649998.5},
    @{Sheet="ALC"; Cell="M112"; Value=-2585.75},
    @{Sheet="ALC"; Cell="N112"; Value=-1652214.5},
    @{Sheet="ALC"; Cell="H137"; Value=822089.6},
    @{Sheet="ALC"; Cell="I137"; Value=1316383.2},
    @{Sheet="ALC"; Cell="K137"; Value=3949149.6},
    @{Sheet="ALC"; Cell="M137"; Value=-3946599.6},
    @{Sheet="ALC"; Cell="H138"; Value=4911.5146},
    @{Sheet="ALC"; Cell="I138"; Value=1548.875},
    @{Sheet="ALC"; Cell="K138"; Value=4646.625},
    @{Sheet="ALC"; Cell="M138"; Value=493.375},
    # --- ARM ---
    @{Sheet="ARM"; Cell="H45"; Value=20845.7},
    @{Sheet="ARM"; Cell="I45"; Value=32998.6},
    @{Sheet="ARM"; Cell="K45"; Value=32998.6},
    @{Sheet="ARM"; Cell="M45"; Value=-32621.6},
    @{Sheet="ARM"; Cell="H61"; Value=6561.15},
    @{Sheet="ARM"; Cell="J61"; Value=4708},
    @{Sheet="ARM"; Cell="L61"; Value=4708},
    @{Sheet="ARM"; Cell="N61"; Value=-5132},
    @{Sheet="ARM"; Cell="H97"; Value=21980.53},
    @{Sheet="ARM"; Cell="I97"; Value=12570.2},
    @{Sheet="ARM"; Cell="J97"; Value=35423.855},
    @{Sheet="ARM"; Cell="K97"; Value=12570.2},
    @{Sheet="ARM"; Cell="L97"; Value=35423.855},
    @{Sheet="ARM"; Cell="M97"; Value=-12074.2},
    @{Sheet="ARM"; Cell="N97"; Value=-36415.855},
    @{Sheet="ARM"; Cell="H102"; Value=4431.3335},
    @{Sheet="ARM"; Cell="I102"; Value=3711},
    @{Sheet="ARM"; Cell="J102"; Value=6592.3335},
    @{Sheet="ARM"; Cell="K102"; Value=3711},
    @{Sheet="ARM"; Cell="L102"; Value=6592.3335},
    @{Sheet="ARM"; Cell="M102"; Value=-2089},
    @{Sheet="ARM"; Cell="N102"; Value=-9836.333500000001},
    @{Sheet="ARM"; Cell="H122"; Value=1002996.06},
    @{Sheet="ARM"; Cell="I122"; Value=2663.6667},
    @{Sheet="ARM"; Cell="K122"; Value=7991.000100000001},
    @{Sheet="ARM"; Cell="M122"; Value=-5541.000100000001},
    @{Sheet="ARM"; Cell="H128"; Value=87500},
    @{Sheet="ARM"; Cell="J128"; Value=87500},
    @{Sheet="ARM"; Cell="L128"; Value=87500},
    @{Sheet="ARM"; Cell="N128"; Value=-97460},
    @{Sheet="ARM"; Cell="H136"; Value=6561.15},
    @{Sheet="ARM"; Cell="J136"; Value=4708},
    @{Sheet="ARM"; Cell="L136"; Value=14124},
    @{Sheet="ARM"; Cell="N136"; Value=-19224},
    # --- BSM ---
    @{Sheet="BSM"; Cell="H107"; Value=3606.52},
    @{Sheet="BSM"; Cell="I107"; Value=3484.476},
    @{Sheet="BSM"; Cell="K107"; Value=3484.476},
    @{Sheet="BSM"; Cell="M107"; Value=-1564.476},
    # --- CRP ---
    @{Sheet="CRP"; Cell="H20"; Value=0},
    @{Sheet="CRP"; Cell="J20"; Value=0},
    @{Sheet="CRP"; Cell="L20"; Value=0},
    @{Sheet="CRP"; Cell="N20"; Value=$null},
    @{Sheet="CRP"; Cell="H22"; Value=515.2222},
    @{Sheet="CRP"; Cell="I22"; Value=569.5},
    @{Sheet="CRP"; Cell="K22"; Value=569.5},
    @{Sheet="CRP"; Cell="M22"; Value=-219.5},
    @{Sheet="CRP"; Cell="H30"; Value=0},
    @{Sheet="CRP"; Cell="J30"; Value=0},
    @{Sheet="CRP"; Cell="L30"; Value=0},
    @{Sheet="CRP"; Cell="N30"; Value=$null},
    @{Sheet="CRP"; Cell="H31"; Value=5438.923},
    @{Sheet="CRP"; Cell="I31"; Value=3115.2856},
    @{Sheet="CRP"; Cell="J31"; Value=8149.8335},
    @{Sheet="CRP"; Cell="K31"; Value=3115.2856},
    @{Sheet="CRP"; Cell="L31"; Value=8149.8335},
    @{Sheet="CRP"; Cell="M31"; Value=-2820.2856},
    @{Sheet="CRP"; Cell="N31"; Value=-8739.833500000001},
    @{Sheet="CRP"; Cell="H34"; Value=5438.923},
    @{Sheet="CRP"; Cell="I34"; Value=3115.2856},
    @{Sheet="CRP"; Cell="J34"; Value=8149.8335},
    @{Sheet="CRP"; Cell="K34"; Value=3115.2856},
    @{Sheet="CRP"; Cell="L34"; Value=8149.8335},
    @{Sheet="CRP"; Cell="M34"; Value=-2913.2856},
    @{Sheet="CRP"; Cell="N34"; Value=-8553.833500000001},
    @{Sheet="CRP"; Cell="H58"; Value=4961.595},
    @{Sheet="CRP"; Cell="I58"; Value=7092.048},
    @{Sheet="CRP"; Cell="K58"; Value=7092.048},
    @{Sheet="CRP"; Cell="M58"; Value=-6889.048},
    @{Sheet="CRP"; Cell="H86"; Value=10347.947},
    @{Sheet="CRP"; Cell="I86"; Value=9817.333000000001},
    @{Sheet="CRP"; Cell="K86"; Value=9817.333000000001},
    @{Sheet="CRP"; Cell="M86"; Value=-8694.333000000001},
    @{Sheet="CRP"; Cell="H89"; Value=10347.947},
    @{Sheet="CRP"; Cell="I89"; Value=9817.333000000001},
    @{Sheet="CRP"; Cell="K89"; Value=49086.665},
    @{Sheet="CRP"; Cell="M89"; Value=-43470.665},
    @{Sheet="CRP"; Cell="H99"; Value=242399.47},
    @{Sheet="CRP"; Cell="J99"; Value=4317.5454},
    @{Sheet="CRP"; Cell="L99"; Value=4317.5454},
    @{Sheet="CRP"; Cell="N99"; Value=-7313.5454},
    @{Sheet="CRP"; Cell="H122"; Value=2589},
    @{Sheet="CRP"; Cell="I122"; Value=2589},
    @{Sheet="CRP"; Cell="J122"; Value=0},
    @{Sheet="CRP"; Cell="K122"; Value=7767},
    @{Sheet="CRP"; Cell="L122"; Value=0},
    @{Sheet="CRP"; Cell="M122"; Value=-5317},
    @{Sheet="CRP"; Cell="N122"; Value=$null},
    @{Sheet="CRP"; Cell="H126"; Value=242399.47},
    @{Sheet="CRP"; Cell="J126"; Value=4317.5454},
    @{Sheet="CRP"; Cell="L126"; Value=12952.6362},
    @{Sheet="CRP"; Cell="N126"; Value=-17892.6362},
    @{Sheet="CRP"; Cell="H127"; Value=55000},
    @{Sheet="CRP"; Cell="J127"; Value=55000},
    @{Sheet="CRP"; Cell="L127"; Value=55000},
    @{Sheet="CRP"; Cell="N127"; Value=-64920},
    @{Sheet="CRP"; Cell="H128"; Value=0},
    @{Sheet="CRP"; Cell="J128"; Value=0},
    @{Sheet="CRP"; Cell="L128"; Value=0},
    @{Sheet="CRP"; Cell="N128"; Value=$null},
    @{Sheet="CRP"; Cell="H132"; Value=26178.908},
    @{Sheet="CRP"; Cell="I132"; Value=8944},
    @{Sheet="CRP"; Cell="J132"; Value=135333.33},
    @{Sheet="CRP"; Cell="K132"; Value=26832},
    @{Sheet="CRP"; Cell="L132"; Value=405999.99},
    @{Sheet="CRP"; Cell="M132"; Value=-24302},
    @{Sheet="CRP"; Cell="N132"; Value=-411059.99},
    @{Sheet="CRP"; Cell="H134"; Value=1843840},
    @{Sheet="CRP"; Cell="I134"; Value=2238133.8},
    @{Sheet="CRP"; Cell="J134"; Value=3802.1667},
    @{Sheet="CRP"; Cell="K134"; Value=6714401.399999999},
    @{Sheet="CRP"; Cell="L134"; Value=11406.5001},
    @{Sheet="CRP"; Cell="M134"; Value=-6711866.399999999},
    @{Sheet="CRP"; Cell="N134"; Value=-16476.5001},
    @{Sheet="CRP"; Cell="H136"; Value=4961.595},
    @{Sheet="CRP"; Cell="I136"; Value=7092.048},
    @{Sheet="CRP"; Cell="K136"; Value=21276.144},
    @{Sheet="CRP"; Cell="M136"; Value=-18726.144},
    # --- CUL ---
    @{Sheet="CUL"; Cell="H131"; Value=10418368},
    @{Sheet="CUL"; Cell="J131"; Value=1720.9438},
    @{Sheet="CUL"; Cell="L131"; Value=5162.8314},
    @{Sheet="CUL"; Cell="N131"; Value=-15242.8314},
    # --- GSM ---
    @{Sheet="GSM"; Cell="H23"; Value=5014},
    @{Sheet="GSM"; Cell="J23"; Value=5014},
    @{Sheet="GSM"; Cell="L23"; Value=5014},
    @{Sheet="GSM"; Cell="M23"; Value=-5460},
    @{Sheet="GSM"; Cell="H97"; Value=5017.2354},
    @{Sheet="GSM"; Cell="I97"; Value=6783.788},
    @{Sheet="GSM"; Cell="J97"; Value=1778.5555},
    @{Sheet="GSM"; Cell="K97"; Value=6783.788},
    @{Sheet="GSM"; Cell="L97"; Value=1778.5555},
    @{Sheet="GSM"; Cell="M97"; Value=-6287.788},
    @{Sheet="GSM"; Cell="N97"; Value=-2770.5555},
    # --- LTW ---
    @{Sheet="LTW"; Cell="H100"; Value=4917.1904},
    @{Sheet="LTW"; Cell="I100"; Value=5133},
    @{Sheet="LTW"; Cell="J100"; Value=4000},
    @{Sheet="LTW"; Cell="K100"; Value=5133},
    @{Sheet="LTW"; Cell="L100"; Value=4000},
    @{Sheet="LTW"; Cell="M100"; Value=-4592},
    @{Sheet="LTW"; Cell="N100"; Value=-5082},
    # --- WVR ---
    @{Sheet="WVR"; Cell="H74"; Value=33620},
    @{Sheet="WVR"; Cell="J74"; Value=41890.08},
    @{Sheet="WVR"; Cell="L74"; Value=41890.08},
    @{Sheet="WVR"; Cell="N74"; Value=-43762.08},
    @{Sheet="WVR"; Cell="H77"; Value=33620},
    @{Sheet="WVR"; Cell="J77"; Value=41890.08},
    @{Sheet="WVR"; Cell="L77"; Value=125670.24},
    @{Sheet="WVR"; Cell="N77"; Value=-135030.24},
    @{Sheet="WVR"; Cell="H122"; Value=3294.2942},
    @{Sheet="WVR"; Cell="I122"; Value=3062.6875},
    @{Sheet="WVR"; Cell="K122"; Value=9188.0625},
    @{Sheet="WVR"; Cell="M122"; Value=-6738.0625},
    @{Sheet="WVR"; Cell="H126"; Value=18968.924},
    @{Sheet="WVR"; Cell="I126"; Value=24862.389},
    @{Sheet="WVR"; Cell="J126"; Value=5708.625},
    @{Sheet="WVR"; Cell="K126"; Value=74587.167},
    @{Sheet="WVR"; Cell="L126"; Value=17125.875},
    @{Sheet="WVR"; Cell="M126"; Value=-72117.167},
    @{Sheet="WVR"; Cell="N126"; Value=-22065.875},
    @{Sheet="WVR"; Cell="H130"; Value=0},
    @{Sheet="WVR"; Cell="J130"; Value=0},
    @{Sheet="WVR"; Cell="L130"; Value=0},
    @{Sheet="WVR"; Cell="N130"; Value=$null}
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}
